$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting period dates (Q4 2021 instead of Q2 2021)
$ws.Range("B8").Value = 44378
$ws.Range("C8").Value = 44561

# H8 used to hold the long "note" text with wrap style; it now holds the
# short "Secretaría Administrativa (UPP)" area text with a plain left style.
$ws.Range("H8").Value = "Secretaría Administrativa (UPP)"
$ws.Range("H8").HorizontalAlignment = -4131
$ws.Range("H8").WrapText = $false

# Validation / update dates
$ws.Range("I8").Value = 44571
$ws.Range("J8").Value = 44571

# K8 now holds the new note text (replacing the old one)
$ws.Range("K8").Value = "Las tablas adyacentes, criterios e hipervinculos que se observan vacios, es por que esta institución no genera inventario de bienes inmuebles, durante el periodo."

# Row height grew to fit the longer note text
$ws.Rows.Item(8).RowHeight = 75

# Column K got a bit narrower
$ws.Columns.Item(11).ColumnWidth = 35.42578125

# Update the view / selection state
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B8").Select()
